$wb = $excel.ActiveWorkbook

# Replace all occurrences of the status text "Ready for handoff" with
# "In Translation" across every worksheet (Overview, zh-cn, de-de).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# The status columns (Overview!E:F, zh-cn!C, de-de!C) are narrower now that
# the text is shorter, matching what Excel's column AutoFit produced for
# the author after the text change.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C1").ColumnWidth = 12.5

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C1").ColumnWidth = 12.5
